# Apply French translations to "Facilitators guidelines - Surface Tension.docx"
#
# Strategy: for each English phrase, locate it (whole-word, case-sensitive
# match on the *entire* text of its run/paragraph) with Find.Execute and then
# overwrite the matched Range's .Text directly. Writing to .Range.Text (as
# opposed to passing the replacement through Find's ReplaceWith parameter)
# keeps straight apostrophes / accented characters exactly as specified and
# does not invoke AutoCorrect "smart quote" substitution, while still
# preserving the run's existing formatting (rPr).
#
# "General VMC Video Introduction" is replaced before the shorter
# "Video Introduction" phrase it contains, so the whole-word match for the
# short phrase can't also fire inside the longer one.

$d = $word.ActiveDocument

function Replace-AllWholeText {
    param(
        [string]$FindText,
        [string]$ReplaceText
    )
    $range = $d.Content
    $guard = 0
    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
    while ($range.Find.Execute($FindText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
        $range.Text = $ReplaceText
        $guard++
        if ($guard -gt 50) { break }
    }
}

Replace-AllWholeText "General VMC Video Introduction" "Vidéo générale introduisant le CVM"

Replace-AllWholeText "Video Title" "Titre de la vidéo"
Replace-AllWholeText "Topic" "Rubrique"
Replace-AllWholeText "Aim(s)" "Objectif(s)"
Replace-AllWholeText "Length" "Durée"
Replace-AllWholeText "Camp Location" "Lieu du camp"
Replace-AllWholeText "Facilitators" "Animateurs"
Replace-AllWholeText "N. of students" "N. des étudiants"
Replace-AllWholeText "Resources" "Les ressources"
Replace-AllWholeText "needed" "nécessaires"
Replace-AllWholeText "Preparations" "Préparations"
Replace-AllWholeText "Video time" "Temps de la vidéo"
Replace-AllWholeText "What facilitator does" "Ce que fait le facilitateur"
Replace-AllWholeText "What learners do" "Ce que font les apprenants"
Replace-AllWholeText "Video Introduction" "Video d'introduction"
Replace-AllWholeText "Introduction of the first experiment" "Introduction de la première expérimentation"
Replace-AllWholeText "Assist the process, provoke thoughts" "Faciliter le processus, susciter des pensées"
